$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.354.93'
$ws.Range('E2').Value = '  +2.48%  '
$ws.Range('D3').Value = '2.106.59'
$ws.Range('E3').Value = '  +0.56%  '
$ws.Range('D4').Value = '''1.007'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '''344.41'
$ws.Range('E5').Value = '  +0.63%  '
$ws.Range('D6').Value = '''1.005'
$ws.Range('E6').Value = '  -0.29%  '
$ws.Range('D7').Value = '''0.5230'
$ws.Range('E7').Value = '  +2.13%  '
$ws.Range('D8').Value = '''0.4446'
$ws.Range('E8').Value = '  +1.13%  '
$ws.Range('D9').Value = '''54.46'
$ws.Range('E9').Value = '  +2.36%  '
$ws.Range('D10').Value = '''0.09444'
$ws.Range('E10').Value = '  +3.57%  '
$ws.Range('E11').Value = '  +0.55%  '
$ws.Range('D12').Value = '''25.04'
$ws.Range('E12').Value = '  +1.73%  '
$ws.Range('D13').Value = '''8.697'
$ws.Range('E13').Value = '  +6.23%  '
$ws.Range('D14').Value = '''6.942'
$ws.Range('E14').Value = '  +2.85%  '
$ws.Range('D15').Value = '2.079.03'
$ws.Range('E15').Value = '  -0.83%  '
$ws.Range('D16').Value = '''101.80'
$ws.Range('E16').Value = '  +2.01%  '
$ws.Range('D17').Value = '''0.00001164'
$ws.Range('E17').Value = '  +1.60%  '
$ws.Range('D18').Value = '''1.007'
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').Value = '''21.24'
$ws.Range('E19').Value = '  +1.13%  '
$ws.Range('D20').Value = '''0.06720'
$ws.Range('E20').Value = '  +1.10%  '
$ws.Range('D21').Value = '''6.375'
$ws.Range('E21').Value = '  +3.28%  '
$ws.Range('E22').Value = '  -0.34%  '
$ws.Range('D23').Value = '30.389.34'
$ws.Range('E23').Value = '  +2.51%  '
$ws.Range('D24').Value = '''12.65'
$ws.Range('E24').Value = '  +0.54%  '
$ws.Range('E25').Value = '  +0.20%  '
$ws.Range('D26').Value = '''22.05'
$ws.Range('E26').Value = '  +1.13%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').Value = '''2.545'
$ws.Range('E27').Value = '  +1.10%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').Value = '''163.27'
$ws.Range('E28').Value = '  +0.56%  '
$ws.Range('D29').Value = '''133.84'
$ws.Range('E29').Value = '  +1.13%  '
$ws.Range('D30').Value = '''1.156'
$ws.Range('E30').Value = '  +2.62%  '
$ws.Range('D31').Value = '''1.739'
$ws.Range('E31').Value = '  +6.28%  '
$ws.Range('D32').Value = '''0.1055'
$ws.Range('E32').Value = '  +1.05%  '
$ws.Range('D33').Value = '''6.854'
$ws.Range('E33').Value = '  +13.41%  '
$ws.Range('D34').Value = '''6.272'
$ws.Range('E34').Value = '  +2.13%  '
$ws.Range('E35').Value = '  -1.10%  '
$ws.Range('D36').Value = '''10.44'
$ws.Range('E36').Value = '  +1.78%  '
$ws.Range('D37').Value = '''0.02629'
$ws.Range('E37').Value = '  +2.62%  '
$ws.Range('D38').Value = '''0.06792'
$ws.Range('E38').Value = '  +2.02%  '
$ws.Range('D39').Value = '''0.7080'
$ws.Range('E39').Value = '  +3.64%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').Value = '''1.350'
$ws.Range('E40').Value = '  +5.10%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').Value = '''12.59'
$ws.Range('E41').Value = '  +2.06%  '
$ws.Range('D42').Value = '''0.2229'
$ws.Range('E42').Value = '  +0.19%  '
$ws.Range('D43').Value = '''0.6859'
$ws.Range('E43').Value = '  +3.34%  '
$ws.Range('D44').Value = '''14.55'
$ws.Range('E44').Value = '  +3.40%  '
$ws.Range('D45').Value = '''2.365'
$ws.Range('E45').Value = '  +3.38%  '
$ws.Range('E46').Value = '  -0.33%  '
$ws.Range('D47').Value = '''1.365'
$ws.Range('E47').Value = '  +17.55%  '
$ws.Range('E48').Value = '  +1.38%  '
$ws.Range('D49').Value = '''0.00000000347'
$ws.Range('E49').Value = '  +4.89%  '
$ws.Range('B50').Value = 'EOS'
$ws.Range('C50').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D50').Value = '''1.220'
$ws.Range('E50').Value = '  +0.15%  '
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').Value = '''1.203'
$ws.Range('E51').Value = '  +9.18%  '
